$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proximity")

$rng = $ws.Range("A5:F6")
$rng.NumberFormat = "@"

$ws.Cells.Item(5, 1).Value = "2026-02-01"
$ws.Cells.Item(5, 2).Value = "15:05:09"
$ws.Cells.Item(5, 3).Value = "15:00"
$ws.Cells.Item(5, 4).Value = "Bedroom Door"
$ws.Cells.Item(5, 5).Value = "EXIT"
$ws.Cells.Item(5, 6).Value = "User EXITED Bedroom"

$ws.Cells.Item(6, 1).Value = "2026-02-01"
$ws.Cells.Item(6, 2).Value = "15:05:26"
$ws.Cells.Item(6, 3).Value = "15:00"
$ws.Cells.Item(6, 4).Value = "Bedroom Door"
$ws.Cells.Item(6, 5).Value = "ENTER"
$ws.Cells.Item(6, 6).Value = "User ENTERED Bedroom"
